$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G3").Value = 'backup@backdoor.com, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg'
$ws.Range("G9").Value = 'Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G10").Value = 'basma.hamed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G12").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G17").Value = 'user@user.com, enas.omran@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G18").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G19").Value = 'neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G20").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G21").Value = 'alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G22").Value = 'backup@backdoor.com, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg'
$ws.Range("G24").Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G28").Value = 'marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg'
$ws.Range("G29").Value = 'esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G31").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G32").Value = 'menna-allah.gamil@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Range("G36").Value = 'user@user.com, enas.omran@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G37").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G38").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G39").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G40").Value = 'alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G41").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G43").Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G47").Value = 'esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg'
$ws.Range("G48").Value = 'maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G49").Value = 'mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range("G50").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G51").Value = 'menna-allah.gamil@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Range("G54").Value = 'Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G56").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G57").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G58").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G59").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G60").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G66").Value = 'Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G67").Value = 'esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G75").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G76").Value = 'neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G77").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G78").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G79").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G81").Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G85").Value = 'Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G86").Value = 'maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G88").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G94").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G95").Value = 'neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G96").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G97").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G98").Value = 'Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G100").Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G104").Value = 'esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg'
$ws.Range("G105").Value = 'basma.hamed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G113").Value = 'ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G115").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
